$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 2023
$ws.Range("B13").Value = 247
$ws.Range("C13").Value = 118
$ws.Range("D13").Value = 1973

$ws.Range("D13").Select()
